# TestTask - Change type ID (int -> Guid)
# Replace integer ID placeholders with GUID strings across the four
# reference sheets (Company, Product, Category, ProductType), format the
# affected cells as Text (and one cell accidentally as Scientific, matching
# the recorded edit), resize columns, fix selections/active sheet and
# tweak the Product sheet's page setup / zoom.

$wb = $excel.ActiveWorkbook

$Company     = $wb.Worksheets.Item(1)
$Product     = $wb.Worksheets.Item(2)
$Category    = $wb.Worksheets.Item(3)
$ProductType = $wb.Worksheets.Item(4)

# ---------------------------------------------------------------------
# 1. Enter the new GUID values. The order below matters: it controls the
#    order in which new entries land in the shared-strings table, so we
#    touch each brand-new GUID for the first time in the same sequence
#    the original authors' edits implied.
# ---------------------------------------------------------------------

$Category.Range("A2").Value    = "4af33f61-8fe2-461b-8eae-cc8344feebe8"
$ProductType.Range("C5").Value = "ff1c323c-123b-4eb4-b3cd-1884bd053b07"
$ProductType.Range("A2").Value = "54d22ad6-5748-4ea7-b7e9-c7a4e0b52220"
$ProductType.Range("A4").Value = "6ae9a401-0a41-4384-8f36-4b67df9846d1"
$ProductType.Range("A6").Value = "36e632d2-98b2-4a1b-8c8f-268aac79271e"
$Company.Range("A2").Value     = "54d32ad6-5748-4ea7-b7e9-c7a4e0b52220"
$Product.Range("A3").Value     = "c0606848-ba9a-41fd-bdf2-d355188803eb"
$Product.Range("A4").Value     = "a00ad315-a4a9-406c-9cb4-b15487b016a9"
$Product.Range("A2").Value     = "c5506848-ba9a-41fd-bdf2-d355188803eb"

# Remaining cells that reuse one of the GUID values above.
$Company.Range("A3").Value  = "4af33f61-8fe2-461b-8eae-cc8344feebe8"

$Product.Range("C2").Value  = "54d32ad6-5748-4ea7-b7e9-c7a4e0b52220"
$Product.Range("E2").Value  = "54d22ad6-5748-4ea7-b7e9-c7a4e0b52220"
$Product.Range("C3").Value  = "54d32ad6-5748-4ea7-b7e9-c7a4e0b52220"
$Product.Range("D3").Value  = "c0606848-ba9a-41fd-bdf2-d355188803eb"
$Product.Range("E3").Value  = "6ae9a401-0a41-4384-8f36-4b67df9846d1"
$Product.Range("C4").Value  = "54d32ad6-5748-4ea7-b7e9-c7a4e0b52220"
$Product.Range("D4").Value  = "c0606848-ba9a-41fd-bdf2-d355188803eb"
$Product.Range("E4").Value  = "36e632d2-98b2-4a1b-8c8f-268aac79271e"

$ProductType.Range("C2").Value = "4af33f61-8fe2-461b-8eae-cc8344feebe8"
$ProductType.Range("C3").Value = "4af33f61-8fe2-461b-8eae-cc8344feebe8"
$ProductType.Range("C6").Value = "ff1c323c-123b-4eb4-b3cd-1884bd053b07"

# ---------------------------------------------------------------------
# 2. Number formats. Apply the (accidental) Scientific format first so it
#    becomes style index 1, then Text ("@") everywhere else so it becomes
#    style index 2 - matching the recorded style table order.
# ---------------------------------------------------------------------

$ProductType.Range("A6").NumberFormat = "0.00E+00"

$Company.Range("A2:A3").NumberFormat = "@"

$Product.Range("A2:E2").NumberFormat = "@"
$Product.Range("A3:E3").NumberFormat = "@"
$Product.Range("A4:E4").NumberFormat = "@"

# ---------------------------------------------------------------------
# 3. Column widths (approximate auto-fit results for the GUID columns).
# ---------------------------------------------------------------------

$Product.Columns.Item(1).ColumnWidth = 43
$Product.Columns.Item(2).ColumnWidth = 36.666666666666664
$Product.Columns.Item(3).ColumnWidth = 43.166666666666664
$Product.Columns.Item(4).ColumnWidth = 39.666666666666664
$Product.Columns.Item(5).ColumnWidth = 37.666666666666664

$Category.Columns.Item(1).ColumnWidth = 42
$Category.Columns.Item(2).ColumnWidth = 27.666666666666668

$ProductType.Columns.Item(1).ColumnWidth = 55.333333333333336
$ProductType.Columns.Item(3).ColumnWidth = 37.333333333333336

# ---------------------------------------------------------------------
# 4. Page setup for the Product sheet.
# ---------------------------------------------------------------------

$Product.PageSetup.PaperSize = 9
$Product.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 5. Selections per sheet (set while each sheet is active so the stored
#    sheetView reflects it), then leave Product as the final active tab.
# ---------------------------------------------------------------------

$Company.Activate()
$Company.Range("A2:A3").Select()

$Category.Activate()
$Category.Range("A2").Select()

$ProductType.Activate()
$ProductType.Range("C1").Select()

$Product.Activate()
$excel.ActiveWindow.Zoom = 115
$Product.Range("D15").Select()
